$d = $word.ActiveDocument

$replacements = @(
    @{old='846÷9=94, 0'; new='495÷7=70, 5'},
    @{old='115÷3=38, 1'; new='368÷8=46, 0'},
    @{old='319÷4=79, 3'; new='775÷7=110, 5'},
    @{old='705÷6=117, 3'; new='807÷7=115, 2'},
    @{old='993÷6=165, 3'; new='781÷8=97, 5'},
    @{old='183÷3=61, 0'; new='780÷8=97, 4'},
    @{old='989÷5=197, 4'; new='673÷4=168, 1'},
    @{old='832÷9=92, 4'; new='512÷9=56, 8'},
    @{old='837÷4=209, 1'; new='675÷2=337, 1'},
    @{old='288÷3=96, 0'; new='759÷2=379, 1'},
    @{old='644÷8=80, 4'; new='211÷9=23, 4'},
    @{old='655÷5=131, 0'; new='502÷8=62, 6'},
    @{old='148÷2=74, 0'; new='821÷8=102, 5'},
    @{old='142÷2=71, 0'; new='172÷2=86, 0'},
    @{old='800÷5=160, 0'; new='784÷7=112, 0'},
    @{old='819÷9=91, 0'; new='251÷5=50, 1'},
    @{old='944÷4=236, 0'; new='811÷4=202, 3'},
    @{old='556÷4=139, 0'; new='953÷9=105, 8'},
    @{old='963÷7=137, 4'; new='623÷7=89, 0'},
    @{old='262÷5=52, 2'; new='350÷7=50, 0'},
    @{old='881÷8=110, 1'; new='812÷7=116, 0'},
    @{old='109÷2=54, 1'; new='245÷4=61, 1'},
    @{old='711÷5=142, 1'; new='703÷7=100, 3'},
    @{old='576÷7=82, 2'; new='622÷5=124, 2'},
    @{old='846÷7=120, 6'; new='480÷6=80, 0'}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

Write-Host "Done applying replacements"
